$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cntf"
$ws.Range("C2").Value = "Il6ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5748063333333333
$ws.Range("H2").Value = 1.724419
$ws.Range("I2").Value = 0.1785342934984892
$ws.Range("J2").Value = 0.1785342934984892
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.338276
$ws.Range("N2").Value = 40.01482799999999
$ws.Range("O2").Value = 0.80809692568033
$ws.Range("P2").Value = 0.8080969256803301
$ws.Range("Q2").Value = 7.666925520547998
$ws.Range("R2").Value = 69.002329684932
$ws.Range("S2").Value = 0.1442730137046388
$ws.Range("T2").Value = 0.1442730137046389

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cntf"
$ws.Range("C3").Value = "Il6ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5748063333333333
$ws.Range("H3").Value = 1.724419
$ws.Range("I3").Value = 0.1785342934984892
$ws.Range("J3").Value = 0.1785342934984892
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.961838
$ws.Range("N3").Value = 8.885514
$ws.Range("O3").Value = 0.1794423943666466
$ws.Range("P3").Value = 0.1794423943666466
$ws.Range("Q3").Value = 1.702483240707333
$ws.Range("R3").Value = 15.322349166366
$ws.Range("S3").Value = 0.03203662110192652
$ws.Range("T3").Value = 0.03203662110192654

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cntf"
$ws.Range("C4").Value = "Il6ra"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5748063333333333
$ws.Range("H4").Value = 1.724419
$ws.Range("I4").Value = 0.1785342934984892
$ws.Range("J4").Value = 0.1785342934984892
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2056733333333333
$ws.Range("N4").Value = 0.6170199999999999
$ws.Range("O4").Value = 0.01246067995302334
$ws.Range("P4").Value = 0.01246067995302335
$ws.Range("Q4").Value = 0.1182223345977777
$ws.Range("R4").Value = 1.06400101138
$ws.Range("S4").Value = 0.00222465869192381
$ws.Range("T4").Value = 0.002224658691923811

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cntf"
$ws.Range("C5").Value = "Il6ra"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.7888463333333333
$ws.Range("H5").Value = 2.366539
$ws.Range("I5").Value = 0.2450149113420932
$ws.Range("J5").Value = 0.2450149113420933
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.338276
$ws.Range("N5").Value = 40.01482799999999
$ws.Range("O5").Value = 0.80809692568033
$ws.Range("P5").Value = 0.8080969256803301
$ws.Range("Q5").Value = 10.521850115588
$ws.Range("R5").Value = 94.69665104029198
$ws.Range("S5").Value = 0.1979957966013841
$ws.Range("T5").Value = 0.1979957966013842

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cntf"
$ws.Range("C6").Value = "Il6ra"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.7888463333333333
$ws.Range("H6").Value = 2.366539
$ws.Range("I6").Value = 0.2450149113420932
$ws.Range("J6").Value = 0.2450149113420933
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.961838
$ws.Range("N6").Value = 8.885514
$ws.Range("O6").Value = 0.1794423943666466
$ws.Range("P6").Value = 0.1794423943666466
$ws.Range("Q6").Value = 2.336435046227333
$ws.Range("R6").Value = 21.027915416046
$ws.Range("S6").Value = 0.04396606234675685
$ws.Range("T6").Value = 0.04396606234675686

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cntf"
$ws.Range("C7").Value = "Il6ra"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.7888463333333333
$ws.Range("H7").Value = 2.366539
$ws.Range("I7").Value = 0.2450149113420932
$ws.Range("J7").Value = 0.2450149113420933
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2056733333333333
$ws.Range("N7").Value = 0.6170199999999999
$ws.Range("O7").Value = 0.01246067995302334
$ws.Range("P7").Value = 0.01246067995302335
$ws.Range("Q7").Value = 0.1622446548644444
$ws.Range("R7").Value = 1.46020189378
$ws.Range("S7").Value = 0.003053052393952213
$ws.Range("T7").Value = 0.003053052393952214

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cntf"
$ws.Range("C8").Value = "Il6ra"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.855932333333334
$ws.Range("H8").Value = 5.567797000000001
$ws.Range("I8").Value = 0.5764507951594176
$ws.Range("J8").Value = 0.5764507951594177
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 13.338276
$ws.Range("N8").Value = 40.01482799999999
$ws.Range("O8").Value = 0.80809692568033
$ws.Range("P8").Value = 0.8080969256803301
$ws.Range("Q8").Value = 24.754937699324
$ws.Range("R8").Value = 222.794439293916
$ws.Range("S8").Value = 0.465828115374307
$ws.Range("T8").Value = 0.4658281153743071

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cntf"
$ws.Range("C9").Value = "Il6ra"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.855932333333334
$ws.Range("H9").Value = 5.567797000000001
$ws.Range("I9").Value = 0.5764507951594176
$ws.Range("J9").Value = 0.5764507951594177
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.961838
$ws.Range("N9").Value = 8.885514
$ws.Range("O9").Value = 0.1794423943666466
$ws.Range("P9").Value = 0.1794423943666466
$ws.Range("Q9").Value = 5.496970910295334
$ws.Range("R9").Value = 49.47273819265801
$ws.Range("S9").Value = 0.1034397109179632
$ws.Range("T9").Value = 0.1034397109179633

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cntf"
$ws.Range("C10").Value = "Il6ra"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.855932333333334
$ws.Range("H10").Value = 5.567797000000001
$ws.Range("I10").Value = 0.5764507951594176
$ws.Range("J10").Value = 0.5764507951594177
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.2056733333333333
$ws.Range("N10").Value = 0.6170199999999999
$ws.Range("O10").Value = 0.01246067995302334
$ws.Range("P10").Value = 0.01246067995302335
$ws.Range("Q10").Value = 0.3817157894377777
$ws.Range("R10").Value = 3.43544210494
$ws.Range("S10").Value = 0.00718296886714732
$ws.Range("T10").Value = 0.007182968867147323
